# Update countries & provincias Spain
# - Refreshes the per-country COVID stats for the countries whose figures
#   changed in this update cycle.
# - Re-sorts the country table (A4:H219) by "Casos totales" (column B)
#   descending, since several countries overtook/were overtaken by
#   neighbours in the ranking once the new totals were applied.
# - Bumps the "Datos actualizados" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A4:A219")

function Set-CountryStats {
    param(
        [string]$Country,
        [double]$CasosTotales,
        [double]$NuevosCasos,
        [double]$CasosActivos,
        [double]$Recuperados,
        [double]$CasosCriticos,
        [double]$MuertesHoy,
        [double]$Muertes
    )

    $cell = $dataRange.Find($Country, $null, $null, 1)
    $r = $cell.Row

    $ws.Cells.Item($r, 2).Value = $CasosTotales
    $ws.Cells.Item($r, 3).Value = $NuevosCasos
    $ws.Cells.Item($r, 4).Value = $CasosActivos
    $ws.Cells.Item($r, 5).Value = $Recuperados
    $ws.Cells.Item($r, 6).Value = $CasosCriticos
    $ws.Cells.Item($r, 7).Value = $MuertesHoy
    $ws.Cells.Item($r, 8).Value = $Muertes
}

Set-CountryStats "Brasil"                332382 1492 135430 175836 0 68 21116
Set-CountryStats "China"                  82971    0  78258     79 0  0  4634
Set-CountryStats "Corea del Sur"          11165   23  10194    705 0  2   266
Set-CountryStats "Bolivia"                 5579  392    575   4774 0 15   230
Set-CountryStats "Guatemala"               2743  231    222   2470 0  3    51
Set-CountryStats "Tayikistan"              2551    0   1089   1418 0  0    44
Set-CountryStats "Haiti"                    812   78     22    765 0  0    25
Set-CountryStats "Principado de Andorra"    762    0    652     59 0  0    51
Set-CountryStats "Uruguay"                  753    0    603    130 0  0    20
Set-CountryStats "Jamaica"                  544   10    191    344 0  0     9

# Re-sort the table descending by "Casos totales" (column B) now that the
# updated figures may have changed the ranking.
$sortRange = $ws.Range("A4:H219")
$sortRange.Sort($ws.Range("B4:B219"), 2)

# Update the "last refreshed" timestamp shown above the table.
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 03:35"
